# Update "想去人数" (want-to-go count) figures in the F column for the
# exhibition list rows on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F value mapping (row 1 is the header row).
$updates = @{
    3 = 2060
    4 = 268
    5 = 70
    6 = 6369
    7 = 254
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
